# Rollerball System Test Cases.xlsx - apply commit changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "Test Steps" text for rows 2 & 3 (register / log in with concrete test data) ---
$ws.Range("C2").Value = "1: Connect to the system`n2: Enter new name(meetsavaliya2), email(meet2@ymail.com), and password(meet123)                       `n3: Click register"
$ws.Range("C3").Value = "1: Connect to the system`n2: Enter existing email(meetsavaliya2) and password(meet123)`n3: Cick log in"

# --- New "Actual results" (E) / "Pass/Fail" (F) columns ---
$ws.Range("E2").Value  = "Done"
$ws.Range("F2").Value  = "Pass"
$ws.Range("E3").Value  = "Done"
$ws.Range("F3").Value  = "Pass"
$ws.Range("E4").Value  = "Done"
$ws.Range("F4").Value  = "Pass"
$ws.Range("E5").Value  = "Done"
$ws.Range("F5").Value  = "Rework"
$ws.Range("E6").Value  = "Done"
$ws.Range("F6").Value  = "Pass"
$ws.Range("E7").Value  = "Done"
$ws.Range("F7").Value  = "Pass"
$ws.Range("E8").Value  = "Done"
$ws.Range("F8").Value  = "Pass"
$ws.Range("E9").Value  = "Done"
$ws.Range("F9").Value  = "Pass"
$ws.Range("E10").Value = "Rework"
$ws.Range("F10").Value = "Rework"
$ws.Range("E11").Value = "Done"
$ws.Range("F11").Value = "Pass"
$ws.Range("E12").Value = "Rework"
$ws.Range("F12").Value = "Rework"

# --- Formatting: wrap text + Arial 10pt font on the new/changed cells ---
$fmtRange = $ws.Range("C2:C4,C7,E2:F12")
$fmtRange.WrapText = $true
$fmtRange.Font.Name = "Arial"
$fmtRange.Font.Size = 10
$ws.Range("E2:F12").Font.Color = 0

# --- Column C widened to fit the longer test-step text ---
$ws.Columns.Item(3).ColumnWidth = 54.15

# --- Row heights recalculated after the content/column changes ---
$ws.Rows.Item(2).RowHeight  = 56
$ws.Rows.Item(3).RowHeight  = 42
$ws.Rows.Item(4).RowHeight  = 42
$ws.Rows.Item(5).RowHeight  = 70
$ws.Rows.Item(6).RowHeight  = 56
$ws.Rows.Item(7).RowHeight  = 70
$ws.Rows.Item(8).RowHeight  = 70
$ws.Rows.Item(9).RowHeight  = 94
$ws.Rows.Item(10).RowHeight = 56
$ws.Rows.Item(11).RowHeight = 70
$ws.Rows.Item(12).RowHeight = 56

# --- View: zoom + selection as last saved ---
$excel.ActiveWindow.Zoom = 131
$ws.Range("E14").Select()
